$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds in rows 4, 6, 7 (odds refreshed since last snapshot) ---
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.25
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("X4").Value = 8.5
$ws.Range("Z4").Value = 13
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 7.5
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 26
$ws.Range("AJ4").Value = 51
$ws.Range("AN4").Value = 8.5
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 41
$ws.Range("AU4").Value = 6.5
$ws.Range("AY4").Value = 101
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3.2
$ws.Range("K7").Value = 1.95
$ws.Range("Q7").Value = 2.4
$ws.Range("R7").Value = 1.53
$ws.Range("S7").Value = 1.53
$ws.Range("T7").Value = 2.38
$ws.Range("U7").Value = 2.05
$ws.Range("V7").Value = 1.7
$ws.Range("Z7").Value = 26
$ws.Range("AC7").Value = 7
$ws.Range("AD7").Value = 6
$ws.Range("AH7").Value = 12
$ws.Range("AR7").Value = 2.38
$ws.Range("BA7").Value = 251
$ws.Range("BD7").Value = 501

# --- Insert a new match row at row 9 (Colombia - Primera A: Santa Fe vs Millonarios) ---
# This pushes the former row 9 (Uruguay match) down to row 10.
$ws.Range("A9").EntireRow.Insert()

$ws.Range("A9").Value = "ph8fDbM8"
$ws.Range("B9").Value = "26/11/2024"
$ws.Range("C9").Value = "22:00"
$ws.Range("D9").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E9").Value = "Santa Fe"
$ws.Range("F9").Value = "Millonarios"
$ws.Range("G9").Value = 3.2
$ws.Range("H9").Value = 2.88
$ws.Range("I9").Value = 2.5
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 1.91
$ws.Range("L9").Value = 3.4
$ws.Range("M9").Value = 1.13
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 1.53
$ws.Range("P9").Value = 2.38
$ws.Range("Q9").Value = 2.7
$ws.Range("R9").Value = 1.44
$ws.Range("S9").Value = 1.62
$ws.Range("T9").Value = 2.2
$ws.Range("U9").Value = 2.2
$ws.Range("V9").Value = 1.62
$ws.Range("W9").Value = 7.5
$ws.Range("X9").Value = 13
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 34
$ws.Range("AA9").Value = 34
$ws.Range("AB9").Value = 41
$ws.Range("AC9").Value = 6
$ws.Range("AD9").Value = 5.5
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 67
$ws.Range("AG9").Value = 6
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 11
$ws.Range("AJ9").Value = 23
$ws.Range("AK9").Value = 26
$ws.Range("AL9").Value = 41
$ws.Range("AM9").Value = 4.75
$ws.Range("AN9").Value = 19
$ws.Range("AO9").Value = 34
$ws.Range("AP9").Value = 67
$ws.Range("AQ9").Value = 101
$ws.Range("AR9").Value = 2.2
$ws.Range("AS9").Value = 9.5
$ws.Range("AT9").Value = 81
$ws.Range("AU9").Value = 4.33
$ws.Range("AV9").Value = 15
$ws.Range("AW9").Value = 29
$ws.Range("AX9").Value = 51
$ws.Range("AY9").Value = 101
$ws.Range("AZ9").Value = 301
$ws.Range("BA9").Value = 351
$ws.Range("BB9").Value = 126
$ws.Range("BC9").Value = 126

# BD9 has no odds for this match (source feed left the column blank), same as the
# empty placeholder cells already present elsewhere in the sheet (e.g. BB4/BC4).
# Force an empty text cell (instead of a truly blank one) and strip the stray
# quote-prefix formatting that typing a lone "'" would otherwise leave behind.
$ws.Range("BD9").Value = "'"
$ws.Range("BC9").Copy()
$ws.Range("BD9").PasteSpecial(-4122)
